$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. PEPList: add a new preset row for the Video Playlist list type
# ------------------------------------------------------------------
$pep = $wb.Worksheets.Item("PEPList")
$pep.Activate()
$pep.Range("A7").Value = "List_VideoPalylist"
$pep.Range("B7").Value = "Video Playlist"
$pep.Range("E7").Value = "click"
$pep.Range("F7").Value = "Fixed list"
$pep.Range("I7").Value = "FixedList(9,9)"
$pep.Range("B7").Select()

# ------------------------------------------------------------------
# 2. FixedList: add the new Video full-list source page
# ------------------------------------------------------------------
$fixedList = $wb.Worksheets.Item("FixedList")
$fixedList.Activate()
$fixedList.Range("A10").Value = "/content/pathology-education/language-masters/en/videos/video-full"
$fixedList.Range("A10").Select()

# ------------------------------------------------------------------
# 3. New "Video" component sheet, placed right after "Image"
# ------------------------------------------------------------------
$imageSheet = $wb.Worksheets.Item("Image")
$video = $wb.Worksheets.Add($null, $imageSheet)
$video.Name = "Video"
$video.Activate()

$video.Range("A1").Value = "TestName"
$video.Range("B1").Value = "dropVideo"
$video.Range("C1").Value = "dropThumbnail"
$video.Range("A2").Value = "Video_Test"
$video.Range("B2").Value = "/content/dam/pathology-education"
$video.Range("C2").Value = "/content/dam/pathology-education"
$video.Range("B2").Select()

# ------------------------------------------------------------------
# 4. TestCases_Final: add the PEPList summary row
# ------------------------------------------------------------------
$testCasesFinal = $wb.Worksheets.Item("TestCases_Final")
$testCasesFinal.Activate()
$testCasesFinal.Range("A9").Value = "PEPList"
$testCasesFinal.Range("B9").Value = "1-6"
$testCasesFinal.Range("B9").NumberFormat = $testCasesFinal.Range("B8").NumberFormat
$testCasesFinal.Range("A2:XFD9").Select()

# ------------------------------------------------------------------
# 5. TestCases: correct the PEPList range and add the Video summary row
# ------------------------------------------------------------------
$testCases = $wb.Worksheets.Item("TestCases")
$testCases.Activate()
$testCases.Range("B9").Value = "1-6"
$testCases.Range("A10").Value = "Video"
$testCases.Range("B10").NumberFormat = $testCases.Range("B8").NumberFormat
$testCases.Range("B10").Value = "1"
$testCases.Range("A2:XFD9").Select()
